$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the data region (columns D, K:T) for rows 2-19 before overwriting,
# since the edit is a permutation of whole rows' data among each other.
$colsD = $ws.Range("D2:D19").Value2
$colsKT = $ws.Range("K2:T19").Value2

# Mapping: new row r gets the data that currently sits in row mapping[r]
# (1-based offsets into the arrays above, where offset 1 = sheet row 2)
$mapping = @{
    2  = 5
    3  = 8
    4  = 9
    5  = 14
    6  = 2
    7  = 3
    8  = 11
    9  = 15
    10 = 18
    11 = 19
    12 = 16
    13 = 17
    14 = 10
    15 = 4
    16 = 12
    17 = 13
    18 = 6
    19 = 7
}

foreach ($destRow in 2..19) {
    $srcRow = $mapping[$destRow]
    $srcOffset = $srcRow - 1   # row index into the snapshot arrays (1-based)

    $ws.Range("D$destRow").Value2 = $colsD[$srcOffset, 1]

    $ws.Range("K$destRow").Value2 = $colsKT[$srcOffset, 1]
    $ws.Range("L$destRow").Value2 = $colsKT[$srcOffset, 2]
    $ws.Range("M$destRow").Value2 = $colsKT[$srcOffset, 3]
    $ws.Range("N$destRow").Value2 = $colsKT[$srcOffset, 4]
    $ws.Range("O$destRow").Value2 = $colsKT[$srcOffset, 5]
    $ws.Range("P$destRow").Value2 = $colsKT[$srcOffset, 6]
    $ws.Range("Q$destRow").Value2 = $colsKT[$srcOffset, 7]
    $ws.Range("R$destRow").Value2 = $colsKT[$srcOffset, 8]
    $ws.Range("S$destRow").Value2 = $colsKT[$srcOffset, 9]
    $ws.Range("T$destRow").Value2 = $colsKT[$srcOffset, 10]
}
